$wb = $excel.ActiveWorkbook

# --- DataFiles sheet: just move the selection (no content changes here) ---
$ws1 = $wb.Worksheets.Item("DataFiles")
$ws1.Range("D1").Select()

# --- tpDictionary sheet: header rename + bugfix ---
$ws2 = $wb.Worksheets.Item("tpDictionary")
$ws2.Activate()

# Convert xlsHeaders to title case
$ws2.Range("A1").Value = "TargetColumn"
$ws2.Range("B1").Value = "Type"
$ws2.Range("C1").Value = "SourceColumn"
$ws2.Range("D1").Value = "SourceUnit"
$ws2.Range("E1").Value = "Filter"
$ws2.Range("F1").Value = "FilterValue"
$ws2.Range("G1").Value = "Description"

# Bugfix: "studyArm" row was mis-typed as "identifier"; should be "metadata"
$ws2.Range("B4").Value = "metadata"

$ws2.Range("A13").Select()
